$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting existing rows 24:46 down to 25:47
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new record's data
$ws.Range("A24").Value = 2
$ws.Range("B24").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C24").Value = "Coquimbo"
$ws.Range("D24").Value = 44671
$ws.Range("E24").Value = 4
$ws.Range("F24").Value = 100112022
$ws.Range("G24").Value = "Arveja Verde"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 240
$ws.Range("K24").Value = 23000
$ws.Range("L24").Value = 25000
$ws.Range("M24").Value = 24000
$ws.Range("N24").Value = "$/malla 25 kilos"
$ws.Range("O24").Value = "Provincia de Limarí"
$ws.Range("P24").Value = 960
$ws.Range("Q24").Value = 25
$ws.Range("R24").Value = "Hortaliza"
